$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed price feed writes plain numeric-looking strings (e.g. "293.59")
# as literal text, not floating point numbers. Pre-format those destination cells
# as Text so Excel stores the literal digits instead of auto-converting them to a
# binary double (which would introduce float rounding noise like 293.58999999999997).
$textPriceCells = @(
    "D5",
    "D6",
    "D10",
    "D11",
    "D12",
    "D14",
    "D16",
    "D18",
    "D21",
    "D22",
    "D23",
    "D24",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D34",
    "D35",
    "D37",
    "D39",
    "D40",
    "D41",
    "D43",
    "D45",
    "D46",
    "D47",
    "D48",
    "D51"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '39.858.04'
$ws.Range("E2").Value = '  -0.43%  '

# Row 3
$ws.Range("D3").Value = '2.204.90'
$ws.Range("E3").Value = '  -1.32%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '293.59'
$ws.Range("E5").Value = '  -0.29%  '

# Row 6
$ws.Range("D6").Value = '86.78'
$ws.Range("E6").Value = '  +0.33%  '

# Row 7
$ws.Range("E7").Value = '  -1.34%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("E9").Value = '  +0.22%  '

# Row 10
$ws.Range("D10").Value = '0.0773'
$ws.Range("E10").Value = '  -2.32%  '

# Row 11
$ws.Range("D11").Value = '29.71'
$ws.Range("E11").Value = '  -4.54%  '

# Row 12
$ws.Range("D12").Value = '49.18'
$ws.Range("E12").Value = '  +4.48%  '

# Row 13
$ws.Range("E13").Value = '  +2.39%  '

# Row 14
$ws.Range("D14").Value = '6.46'
$ws.Range("E14").Value = '  +0.04%  '

# Row 15
$ws.Range("D15").Value = '2.550.89'
$ws.Range("E15").Value = '  -1.19%  '

# Row 16
$ws.Range("D16").Value = '13.71'
$ws.Range("E16").Value = '  -3.16%  '

# Row 17
$ws.Range("D17").Value = '2.221.06'
$ws.Range("E17").Value = '  -0.57%  '

# Row 18
$ws.Range("D18").Value = '0.726'
$ws.Range("E18").Value = '  -0.35%  '

# Row 19
$ws.Range("D19").Value = '39.790.42'
$ws.Range("E19").Value = '  -0.43%  '

# Row 20
$ws.Range("E20").Value = '  -0.82%  '

# Row 21
$ws.Range("D21").Value = '11.25'
$ws.Range("E21").Value = '  +3.94%  '

# Row 22
$ws.Range("D22").Value = '5.75'
$ws.Range("E22").Value = '  -0.87%  '

# Row 23
$ws.Range("D23").Value = '65.08'
$ws.Range("E23").Value = '  -0.52%  '

# Row 24
$ws.Range("D24").Value = '235.91'
$ws.Range("E24").Value = '  +0.31%  '

# Row 25
$ws.Range("E25").Value = '  +0.02%  '

# Row 26
$ws.Range("D26").Value = '2.45'
$ws.Range("E26").Value = '  -0.69%  '

# Row 27
$ws.Range("E27").Value = '  -2.19%  '

# Row 28
$ws.Range("D28").Value = '22.43'
$ws.Range("E28").Value = '  -1.59%  '

# Row 29
$ws.Range("D29").Value = '2.15'
$ws.Range("E29").Value = '  -4.31%  '

# Row 30
$ws.Range("D30").Value = '9.15'
$ws.Range("E30").Value = '  -0.76%  '

# Row 31
$ws.Range("D31").Value = '155.19'
$ws.Range("E31").Value = '  +1.88%  '

# Row 32
$ws.Range("D32").Value = '31.74'
$ws.Range("E32").Value = '  -4.75%  '

# Row 33
$ws.Range("E33").Value = '  +0.03%  '

# Row 34
$ws.Range("D34").Value = '4.88'
$ws.Range("E34").Value = '  -0.02%  '

# Row 35
$ws.Range("D35").Value = '0.0710'
$ws.Range("E35").Value = '  -0.96%  '

# Row 36
$ws.Range("E36").Value = '  -1.40%  '

# Row 37
$ws.Range("D37").Value = '2.82'
$ws.Range("E37").Value = '  +4.13%  '

# Row 38
$ws.Range("E38").Value = '  -0.18%  '

# Row 39
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '0.0974'
$ws.Range("E39").Value = '  -2.74%  '

# Row 40
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").Value = '15.44'
$ws.Range("E40").Value = '  -5.75%  '

# Row 41
$ws.Range("D41").Value = '1.66'
$ws.Range("E41").Value = '  -2.18%  '

# Row 42
$ws.Range("D42").Value = '2.118.42'
$ws.Range("E42").Value = '  +3.67%  '

# Row 43
$ws.Range("D43").Value = '3.73'
$ws.Range("E43").Value = '  -2.96%  '

# Row 44
$ws.Range("E44").Value = '  -5.87%  '

# Row 45
$ws.Range("D45").Value = '0.0266'
$ws.Range("E45").Value = '  -1.43%  '

# Row 46
$ws.Range("D46").Value = '17.65'
$ws.Range("E46").Value = '  +8.07%  '

# Row 47
$ws.Range("D47").Value = '9.62'
$ws.Range("E47").Value = '  -3.53%  '

# Row 48
$ws.Range("D48").Value = '2.65'
$ws.Range("E48").Value = '  +3.61%  '

# Row 49
$ws.Range("D49").Value = '2.417.21'
$ws.Range("E49").Value = '  -1.43%  '

# Row 50
$ws.Range("E50").Value = '  -0.31%  '

# Row 51
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = '1.09'
$ws.Range("E51").Value = '  -0.21%  '
